$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 71449.07000000001
$ws.Range("I11").Value = 71449.07000000001
$ws.Range("K11").Value = 71449.07000000001
$ws.Range("M11").Value = -71309.07000000001

$ws.Range("H106").Value = 3616.4
$ws.Range("I106").Value = 3573.7778
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 3573.7778
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -2942.7778
$ws.Range("N106").Value = -5262

$ws.Range("H107").Value = 696.1539
$ws.Range("I107").Value = 683.44446
$ws.Range("J107").Value = 724.75
$ws.Range("K107").Value = 683.44446
$ws.Range("L107").Value = 724.75
$ws.Range("M107").Value = 1236.55554
$ws.Range("N107").Value = -4564.75

$ws.Range("H113").Value = 51425.2
$ws.Range("I113").Value = 112356
$ws.Range("K113").Value = 112356
$ws.Range("M113").Value = -109102

$ws.Range("H116").Value = 2249.625
$ws.Range("J116").Value = 1999.25
$ws.Range("L116").Value = 1999.25
$ws.Range("N116").Value = -8883.25

$ws.Range("H138").Value = 3471.75
$ws.Range("I138").Value = 2212.8696
$ws.Range("J138").Value = 4062.653
$ws.Range("K138").Value = 6638.6088
$ws.Range("L138").Value = 12187.959
$ws.Range("M138").Value = -1498.6088
$ws.Range("N138").Value = -22467.959

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 84721.75
$ws.Range("I2").Value = 1566.1
$ws.Range("J2").Value = 500500
$ws.Range("K2").Value = 1566.1
$ws.Range("L2").Value = 500500
$ws.Range("M2").Value = -1453.1
$ws.Range("N2").Value = -500726

$ws.Range("H32").Value = 38395.902
$ws.Range("I32").Value = 16313.033
$ws.Range("J32").Value = 173101.4
$ws.Range("K32").Value = 16313.033
$ws.Range("L32").Value = 173101.4
$ws.Range("M32").Value = -16026.033
$ws.Range("N32").Value = -173675.4

$ws.Range("H44").Value = 12855.571
$ws.Range("J44").Value = 12831.5
$ws.Range("L44").Value = 12831.5
$ws.Range("N44").Value = -13807.5

$ws.Range("H55").Value = 12128.571
$ws.Range("J55").Value = 12128.571
$ws.Range("L55").Value = 12128.571
$ws.Range("N55").Value = -12758.571

$ws.Range("H63").Value = 2600.5715
$ws.Range("I63").Value = 902.5
$ws.Range("K63").Value = 902.5
$ws.Range("M63").Value = -216.5

$ws.Range("H66").Value = 2600.5715
$ws.Range("I66").Value = 902.5
$ws.Range("K66").Value = 4512.5
$ws.Range("M66").Value = -1080.5

$ws.Range("H80").Value = 25744.545
$ws.Range("J80").Value = 25744.545
$ws.Range("L80").Value = 25744.545
$ws.Range("N80").Value = -27740.545

$ws.Range("H83").Value = 25744.545
$ws.Range("J83").Value = 25744.545
$ws.Range("L83").Value = 77233.63499999999
$ws.Range("N83").Value = -87217.63499999999

$ws.Range("H116").Value = 84721.75
$ws.Range("I116").Value = 1566.1
$ws.Range("J116").Value = 500500
$ws.Range("K116").Value = 1566.1
$ws.Range("L116").Value = 500500
$ws.Range("M116").Value = 727.9000000000001
$ws.Range("N116").Value = -505088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 84721.75
$ws.Range("I3").Value = 1566.1
$ws.Range("J3").Value = 500500
$ws.Range("K3").Value = 1566.1
$ws.Range("L3").Value = 500500
$ws.Range("M3").Value = -1452.1
$ws.Range("N3").Value = -500728

$ws.Range("H22").Value = 397.14285
$ws.Range("J22").Value = 397.5
$ws.Range("L22").Value = 397.5
$ws.Range("N22").Value = -743.5

$ws.Range("H80").Value = 1099.5758
$ws.Range("I80").Value = 640.2857
$ws.Range("J80").Value = 1438
$ws.Range("K80").Value = 640.2857
$ws.Range("L80").Value = 1438
$ws.Range("M80").Value = 357.7143
$ws.Range("N80").Value = -3434

$ws.Range("H82").Value = 18975.334
$ws.Range("J82").Value = 34786.668
$ws.Range("L82").Value = 34786.668
$ws.Range("N82").Value = -35552.668

$ws.Range("H83").Value = 1099.5758
$ws.Range("I83").Value = 640.2857
$ws.Range("J83").Value = 1438
$ws.Range("K83").Value = 3201.4285
$ws.Range("L83").Value = 7190
$ws.Range("M83").Value = 1790.5715
$ws.Range("N83").Value = -17174

$ws.Range("H85").Value = 18975.334
$ws.Range("J85").Value = 34786.668
$ws.Range("L85").Value = 34786.668
$ws.Range("N85").Value = -37438.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33091.625
$ws.Range("I31").Value = 1571.25
$ws.Range("J31").Value = 64612
$ws.Range("K31").Value = 1571.25
$ws.Range("L31").Value = 64612
$ws.Range("M31").Value = -1276.25
$ws.Range("N31").Value = -65202

$ws.Range("H34").Value = 33091.625
$ws.Range("I34").Value = 1571.25
$ws.Range("J34").Value = 64612
$ws.Range("K34").Value = 1571.25
$ws.Range("L34").Value = 64612
$ws.Range("M34").Value = -1369.25
$ws.Range("N34").Value = -65016

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 545
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 990
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 2970
$ws.Range("M13").Value = -132
$ws.Range("N13").Value = -3306

$ws.Range("H137").Value = 40360.934
$ws.Range("I137").Value = 86280
$ws.Range("J137").Value = 9748.223
$ws.Range("K137").Value = 258840
$ws.Range("L137").Value = 29244.669
$ws.Range("M137").Value = -253740
$ws.Range("N137").Value = -39444.669

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1263113.9
$ws.Range("I107").Value = 483.5
$ws.Range("J107").Value = 5051005
$ws.Range("K107").Value = 483.5
$ws.Range("L107").Value = 5051005
$ws.Range("M107").Value = 1436.5
$ws.Range("N107").Value = -5054845

$ws.Range("H132").Value = 2754.0889
$ws.Range("I132").Value = 2031.909
$ws.Range("K132").Value = 6095.727000000001
$ws.Range("M132").Value = -3565.727000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 74059.86
$ws.Range("I40").Value = 252669.5
$ws.Range("J40").Value = 2616
$ws.Range("K40").Value = 252669.5
$ws.Range("L40").Value = 2616
$ws.Range("M40").Value = -252533.5
$ws.Range("N40").Value = -2888

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6070
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H81").Value = 200544.3
$ws.Range("I81").Value = 143413.14
$ws.Range("J81").Value = 333850.34
$ws.Range("K81").Value = 286826.28
$ws.Range("L81").Value = 667700.6800000001
$ws.Range("M81").Value = -285765.28
$ws.Range("N81").Value = -669822.6800000001

$ws.Range("H84").Value = 200544.3
$ws.Range("I84").Value = 143413.14
$ws.Range("J84").Value = 333850.34
$ws.Range("K84").Value = 1434131.4
$ws.Range("L84").Value = 3338503.4
$ws.Range("M84").Value = -1428827.4
$ws.Range("N84").Value = -3349111.4

$ws.Range("H119").Value = 29950
$ws.Range("J119").Value = 29950
$ws.Range("L119").Value = 29950
$ws.Range("N119").Value = -39626
